# Full copyedit: drop the separate bold title paragraph ("Mixing cover
# crops suppresses weeds and roto-till improves urban soil compaction and
# infiltration") and promote the following "Highlights:" paragraph (which
# used to carry the BodyText style) to the FirstParagraph style so it now
# opens the document.

$d = $word.ActiveDocument

# 1) Remove the first paragraph (title) entirely, including its paragraph
#    mark, so the "Highlights:" paragraph shifts up to become paragraph 1.
$d.Paragraphs.Item(1).Range.Delete()

# 2) The paragraph that used to be "Highlights:" (style BodyText) is now
#    paragraph 1. Re-style it as FirstParagraph.
$target = $d.Paragraphs.Item(1)
$target.Style = "First Paragraph"

# Changing the paragraph style can reset direct run-level character
# formatting that was already present on the text (bold / complex-script
# bold). Restore it explicitly on the text run, excluding the trailing
# paragraph mark.
$textRange = $d.Range($target.Range.Start, $target.Range.End - 1)
$textRange.Font.Bold = -1
$textRange.Font.BoldBi = -1
